$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-09-24 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-09-25 Thursday", 2)

# Update the division problems in the table. Table cells are addressed by
# (row, column) because several cells share identical text (e.g. "92÷4="),
# so a global Find/Replace would not be safe.
$table = $d.Tables.Item(1)

$table.Cell(1, 1).Range.Text = "21÷8="
$table.Cell(1, 2).Range.Text = "20÷5="
$table.Cell(1, 3).Range.Text = "69÷4="
$table.Cell(1, 4).Range.Text = "81÷8="
$table.Cell(1, 5).Range.Text = "74÷7="

$table.Cell(5, 1).Range.Text = "16÷4="
$table.Cell(5, 2).Range.Text = "87÷4="
$table.Cell(5, 3).Range.Text = "15÷3="
$table.Cell(5, 4).Range.Text = "56÷5="
$table.Cell(5, 5).Range.Text = "69÷7="

$table.Cell(9, 1).Range.Text = "12÷3="
$table.Cell(9, 2).Range.Text = "67÷4="
$table.Cell(9, 3).Range.Text = "53÷8="
$table.Cell(9, 4).Range.Text = "66÷3="
$table.Cell(9, 5).Range.Text = "33÷7="

$table.Cell(13, 1).Range.Text = "25÷5="
$table.Cell(13, 2).Range.Text = "70÷3="
$table.Cell(13, 3).Range.Text = "31÷5="
$table.Cell(13, 4).Range.Text = "38÷8="
$table.Cell(13, 5).Range.Text = "62÷3="

$table.Cell(17, 1).Range.Text = "76÷4="
$table.Cell(17, 2).Range.Text = "91÷2="
$table.Cell(17, 3).Range.Text = "24÷6="
$table.Cell(17, 4).Range.Text = "92÷8="
$table.Cell(17, 5).Range.Text = "84÷6="
